$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook has 3 sheets (left-to-right tab order, which stays fixed):
#   position 1 (rId1): currently named "2021-Q2" - holds the fund-detail table (A1:H24)
#   position 2 (rId2): "2021-Q3" - untouched by this edit
#   position 3 (rId3): currently named "总计"    - holds the summary table (A1:D3)
#
# The edit swaps the *content* (and name) that live at position 1 and
# position 3, so that afterwards:
#   position 1 -> name "总计", holds the summary table
#   position 3 -> name "2021-Q2", holds the fund-detail table (and becomes
#                 the active/selected tab)
# position 2 ("2021-Q3") keeps its name & content unchanged.
#
# NOTE: worksheet object variables captured via Worksheets.Item(...) behave
# like positional handles in this runtime: once the sheet collection is
# restructured (Add/Delete), a previously captured variable can end up
# referring to a *different* sheet. To stay safe, every sheet is re-fetched
# by name immediately before each use instead of being cached in a variable.
# ---------------------------------------------------------------------------

# Create a temporary holding sheet, and give it a stable, unique name right away.
$tmp = $wb.Worksheets.Add()
$tmp.Name = "TempSwapSheet"

# Step 1: move the fund-detail content ("2021-Q2") into the temp sheet.
$wb.Worksheets.Item("2021-Q2").Range("B1:H24").Copy($wb.Worksheets.Item("TempSwapSheet").Range("B1"))
$wb.Worksheets.Item("2021-Q2").Range("A2:A24").Copy($wb.Worksheets.Item("TempSwapSheet").Range("A2"))
$wb.Worksheets.Item("2021-Q2").Cells.Clear()

# Step 2: move the summary content ("总计") into the now-empty "2021-Q2" sheet.
$wb.Worksheets.Item("总计").Range("B1:D3").Copy($wb.Worksheets.Item("2021-Q2").Range("B1"))
$wb.Worksheets.Item("总计").Range("A2:A3").Copy($wb.Worksheets.Item("2021-Q2").Range("A2"))
$wb.Worksheets.Item("总计").Cells.Clear()

# Step 3: move the fund-detail content out of the temp sheet into the now-empty "总计" sheet.
$wb.Worksheets.Item("TempSwapSheet").Range("B1:H24").Copy($wb.Worksheets.Item("总计").Range("B1"))
$wb.Worksheets.Item("TempSwapSheet").Range("A2:A24").Copy($wb.Worksheets.Item("总计").Range("A2"))

# Remove the temp helper sheet.
$wb.Worksheets.Item("TempSwapSheet").Delete()

# Step 4: swap the page margins that go with each table (the summary table's
# sheet used 0.75/1/0.5 in / the fund-detail table's sheet used 0.7/0.75/0.3 in).
$wb.Worksheets.Item("2021-Q2").PageSetup.LeftMargin = 54
$wb.Worksheets.Item("2021-Q2").PageSetup.RightMargin = 54
$wb.Worksheets.Item("2021-Q2").PageSetup.TopMargin = 72
$wb.Worksheets.Item("2021-Q2").PageSetup.BottomMargin = 72
$wb.Worksheets.Item("2021-Q2").PageSetup.HeaderMargin = 36
$wb.Worksheets.Item("2021-Q2").PageSetup.FooterMargin = 36

$wb.Worksheets.Item("总计").PageSetup.LeftMargin = 50.4
$wb.Worksheets.Item("总计").PageSetup.RightMargin = 50.4
$wb.Worksheets.Item("总计").PageSetup.TopMargin = 54
$wb.Worksheets.Item("总计").PageSetup.BottomMargin = 54
$wb.Worksheets.Item("总计").PageSetup.HeaderMargin = 21.6
$wb.Worksheets.Item("总计").PageSetup.FooterMargin = 21.6

# Step 5: rename the sheets, swapping "2021-Q2" <-> "总计" names.
# (use a temporary unique name to avoid a duplicate-name collision mid-swap)
$wb.Worksheets.Item("2021-Q2").Name = "TempRenameXYZ"
$wb.Worksheets.Item("总计").Name = "2021-Q2"
$wb.Worksheets.Item("TempRenameXYZ").Name = "总计"

# Step 6: the sheet that now holds the fund-detail data (named "2021-Q2") becomes
# the selected/active tab, matching the original layout where that table's
# sheet was the active one.
$wb.Worksheets.Item("2021-Q2").Activate()
